$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was updated from 45206 (2023-10-07)
# to 45208 (2023-10-09) for every data row (C2:C399).
$ws.Range("C2:C399").Value = 45208
